$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 20000
$ws.Range("J16").Value = 20000
$ws.Range("L16").Value = 20000
$ws.Range("N16").Value = -20460
$ws.Range("H76").Value = 6149.1
$ws.Range("I76").Value = 8400.6
$ws.Range("J76").Value = 3897.6
$ws.Range("K76").Value = 8400.6
$ws.Range("L76").Value = 3897.6
$ws.Range("M76").Value = -8085.6
$ws.Range("N76").Value = -4527.6
$ws.Range("H79").Value = 6149.1
$ws.Range("I79").Value = 8400.6
$ws.Range("J79").Value = 3897.6
$ws.Range("K79").Value = 8400.6
$ws.Range("L79").Value = 3897.6
$ws.Range("M79").Value = -7308.6
$ws.Range("N79").Value = -6081.6
$ws.Range("H134").Value = 138116.67
$ws.Range("J134").Value = 138116.67
$ws.Range("L134").Value = 138116.67
$ws.Range("N134").Value = -148256.67
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5285
$ws.Range("I132").Value = 5898.5415
$ws.Range("J132").Value = 3648.889
$ws.Range("K132").Value = 17695.6245
$ws.Range("L132").Value = 10946.667
$ws.Range("M132").Value = -15165.6245
$ws.Range("N132").Value = -16006.667
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 8171.3335
$ws.Range("I33").Value = 2000
$ws.Range("J33").Value = 11257
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 11257
$ws.Range("M33").Value = -1664
$ws.Range("N33").Value = -11929
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49999
$ws.Range("J20").Value = 49999
$ws.Range("L20").Value = 49999
$ws.Range("N20").Value = -50471
$ws.Range("H30").Value = 49999
$ws.Range("J30").Value = 49999
$ws.Range("L30").Value = 49999
$ws.Range("N30").Value = -50181
$ws.Range("H128").Value = 49999
$ws.Range("J128").Value = 49999
$ws.Range("L128").Value = 49999
$ws.Range("N128").Value = -59959
$ws.Range("H132").Value = 522056.16
$ws.Range("I132").Value = 712836.75
$ws.Range("J132").Value = 4223.143
$ws.Range("K132").Value = 2138510.25
$ws.Range("L132").Value = 12669.429
$ws.Range("M132").Value = -2135980.25
$ws.Range("N132").Value = -17729.429
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7038.7827
$ws.Range("I3").Value = 3744
$ws.Range("J3").Value = 7954
$ws.Range("K3").Value = 11232
$ws.Range("L3").Value = 23862
$ws.Range("M3").Value = -11120
$ws.Range("N3").Value = -24086
$ws.Range("H5").Value = 2098.5
$ws.Range("I5").Value = 2098.5
$ws.Range("K5").Value = 6295.5
$ws.Range("M5").Value = -6183.5
$ws.Range("H34").Value = 674.3913
$ws.Range("J34").Value = 1190.1
$ws.Range("L34").Value = 3570.3
$ws.Range("N34").Value = -3738.3
$ws.Range("H39").Value = 4547.9
$ws.Range("J39").Value = 5559.875
$ws.Range("L39").Value = 16679.625
$ws.Range("N39").Value = -17267.625
$ws.Range("H55").Value = 6678
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 8300.286
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 24900.858
$ws.Range("M55").Value = -2823
$ws.Range("N55").Value = -25254.858
$ws.Range("H82").Value = 20130
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 20130
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 60390
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -61202
$ws.Range("H85").Value = 20130
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 20130
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 60390
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -63198
$ws.Range("H113").Value = 632.72974
$ws.Range("I113").Value = 591.84
$ws.Range("J113").Value = 717.9167
$ws.Range("K113").Value = 1775.52
$ws.Range("L113").Value = 2153.7501
$ws.Range("M113").Value = 394.48
$ws.Range("N113").Value = -6493.7501
$ws.Range("H121").Value = 111042.6
$ws.Range("I121").Value = 1590
$ws.Range("J121").Value = 184011
$ws.Range("K121").Value = 4770
$ws.Range("L121").Value = 552033
$ws.Range("M121").Value = -3460
$ws.Range("N121").Value = -554653
$ws.Range("H122").Value = 999.6667
$ws.Range("J122").Value = 999.6667
$ws.Range("L122").Value = 8997.0003
$ws.Range("N122").Value = -13897.0003
$ws.Range("H129").Value = 4546234
$ws.Range("I129").Value = 392.5
$ws.Range("J129").Value = 7143858
$ws.Range("K129").Value = 1177.5
$ws.Range("L129").Value = 21431574
$ws.Range("M129").Value = 3822.5
$ws.Range("N129").Value = -21441574
$ws.Range("H134").Value = 3312.55
$ws.Range("I134").Value = 1586.5217
$ws.Range("J134").Value = 5647.7646
$ws.Range("K134").Value = 4759.5651
$ws.Range("L134").Value = 16943.2938
$ws.Range("M134").Value = 310.4349000000002
$ws.Range("N134").Value = -27083.2938
$ws.Range("H135").Value = 2098.5
$ws.Range("I135").Value = 2098.5
$ws.Range("K135").Value = 18886.5
$ws.Range("M135").Value = -16351.5
$ws.Range("H137").Value = 1685.6562
$ws.Range("I137").Value = 877.1429
$ws.Range("J137").Value = 2314.5
$ws.Range("K137").Value = 2631.4287
$ws.Range("L137").Value = 6943.5
$ws.Range("M137").Value = 2468.5713
$ws.Range("N137").Value = -17143.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H126").Value = 3006.3333
$ws.Range("I126").Value = 2881.818
$ws.Range("J126").Value = 3202
$ws.Range("K126").Value = 8645.454000000002
$ws.Range("L126").Value = 9606
$ws.Range("M126").Value = -6175.454000000002
$ws.Range("N126").Value = -14546
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3747.5
$ws.Range("I62").Value = 3330
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3330
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2706
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 3747.5
$ws.Range("I65").Value = 3330
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 16650
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -13530
$ws.Range("N65").Value = -31240
$ws.Range("H132").Value = 2434.6667
$ws.Range("I132").Value = 1841.9474
$ws.Range("J132").Value = 2997.75
$ws.Range("K132").Value = 5525.8422
$ws.Range("L132").Value = 8993.25
$ws.Range("M132").Value = -2995.8422
$ws.Range("N132").Value = -14053.25
